$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "MB Endpoint" column (H) ---------------------------------------
# Values per row (header row 1 + 23 data rows, rows 2-24)
$mbEndpoint = @{
    1  = "MB Endpoint"
    2  = "Yes"
    3  = "Yes"
    4  = "No"
    5  = "Yes"
    6  = "No"
    7  = "No"
    8  = "No"
    9  = "No"
    10 = "No"
    11 = "No"
    12 = "Yes"
    13 = "No"
    14 = "No"
    15 = "No"
    16 = "No"
    17 = "No"
    18 = "No"
    19 = "No"
    20 = "No"
    21 = "No"
    22 = "No"
    23 = "No"
    24 = "No"
}

# Use the plain formatting already applied across column G (row 1's cell,
# which carries the same style as most of the sheet) as the template for
# every new H cell so the new column matches the rest of the table.
$ws.Range("G1").Copy()
foreach ($row in 1..24) {
    $ws.Cells.Item($row, 8).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($row in 1..24) {
    $ws.Cells.Item($row, 8).Value = $mbEndpoint[$row]
}

# --- Existing "N/A" value for NCT05564949 (row 22) is now "None" -------
$ws.Range("G22").Value = "None"

# --- Selection moved to reflect the newly-added column ------------------
[void]$ws.Range("H25").Select()
